# Adding handling for URLs in broader column
$wb = $excel.ActiveWorkbook

# The edit targets the "testreg4" worksheet (second tab / sheet2.xml)
$ws = $wb.Worksheets.Item("testreg4")

# Disable autocorrect of hyperlinks so the URL text in column H (broader)
# stays a plain string rather than being turned into a live hyperlink.
$excel.AutoCorrect.AutoFillHyperlink = $false

# Append a new row of data (row 7), matching columns:
# A=id, B=label, C=altLabel, D=description, E=notation, F=note, G=source, H=broader
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = "item 6"
$ws.Cells.Item(7, 3).Value = "atl label 6"
$ws.Cells.Item(7, 4).Value = "description 6"
$ws.Cells.Item(7, 5).Value = "i6"
$ws.Cells.Item(7, 8).Value = "http://registry.it.csiro.au/sandbox/csiro/utils/commondef/1"

# Remove any hyperlink Excel may have auto-created from the URL text
if ($ws.Hyperlinks.Count -gt 0) {
    $ws.Hyperlinks.Delete()
}

# Update the selection to match the post-edit state recorded in the sheet
$ws.Range("H12").Select()
